# Apply updated dSF (column F) values for the rows whose "repulled" data changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -3
    7  = -3
    13 = -2
    17 = 3
    19 = -1
    20 = 4
    26 = 2
    28 = -2
    29 = -2
    37 = -9
    41 = -4
    43 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
